$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.942.68'
$ws.Range('E2').Value = '  -0.14%  '
$ws.Range('D3').Value = '2.552.88'
$ws.Range('E3').Value = '  +0.38%  '
$ws.Range('E4').Value = '  -0.18%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '303.90'
$ws.Range('E5').Value = '  +1.43%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '98.23'
$ws.Range('E6').Value = '  +6.28%  '
$ws.Range('E7').Value = '  +0.80%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('E9').Value = '  -0.16%  '
$ws.Range('E10').Value = '  +4.30%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0828'
$ws.Range('E11').Value = '  +3.26%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.117'
$ws.Range('E12').Value = '  +5.42%  '
$ws.Range('E13').Value = '  -0.08%  '
$ws.Range('D14').Value = '2.946.18'
$ws.Range('E14').Value = '  +0.31%  '
$ws.Range('D15').Value = '2.597.76'
$ws.Range('E15').Value = '  +0.33%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '15.09'
$ws.Range('E16').Value = '  +7.21%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.878'
$ws.Range('E17').Value = '  +1.21%  '
$ws.Range('D18').Value = '42.980.04'
$ws.Range('E18').Value = '  -0.17%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '13.87'
$ws.Range('E19').Value = '  +6.24%  '
$ws.Range('E20').Value = '  +2.16%  '
$ws.Range('E21').Value = '  +0.25%  '
$ws.Range('E22').Value = '  +0.53%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '254.60'
$ws.Range('E23').Value = '  -0.86%  '
$ws.Range('E24').Value = '  +1.99%  '
$ws.Range('E25').Value = '  -1.76%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '28.05'
$ws.Range('E26').Value = '  -3.24%  '
$ws.Range('E27').Value = '  -0.15%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '10.24'
$ws.Range('E28').Value = '  +2.22%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '37.87'
$ws.Range('E29').Value = '  +1.00%  '
$ws.Range('E30').Value = '  -1.30%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '6.16'
$ws.Range('E31').Value = '  +3.56%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '158.64'
$ws.Range('E32').Value = '  +2.84%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '19.65'
$ws.Range('E33').Value = '  +16.70%  '
$ws.Range('E34').Value = '  +0.36%  '
$ws.Range('B35').Value = 'Hedera'
$ws.Range('C35').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.0804'
$ws.Range('E35').Value = '  +1.37%  '
$ws.Range('B36').Value = 'LidoDAOToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '3.31'
$ws.Range('E36').Value = '  -1.70%  '
$ws.Range('B37').Value = 'WEMIXToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '2.63'
$ws.Range('E37').Value = '  -4.55%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.116'
$ws.Range('E38').Value = '  +1.27%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '25.53'
$ws.Range('E39').Value = '  +9.50%  '
$ws.Range('E40').Value = '  +0.16%  '
$ws.Range('B41').Value = 'ApeXProtocol'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '2.09'
$ws.Range('E41').Value = '  +31.70%  '
$ws.Range('B42').Value = 'NEARProtocol'
$ws.Range('C42').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '3.43'
$ws.Range('E42').Value = '  +0.19%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '3.91'
$ws.Range('E43').Value = '  +0.58%  '
$ws.Range('E44').Value = '  -1.16%  '
$ws.Range('D45').Value = '2.089.41'
$ws.Range('E45').Value = '  +0.61%  '
$ws.Range('E46').Value = '  -0.04%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '86.69'
$ws.Range('E47').Value = '  +2.58%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '8.99'
$ws.Range('E48').Value = '  +1.04%  '
$ws.Range('D49').Value = '2.803.71'
$ws.Range('E49').Value = '  +0.26%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '74.92'
$ws.Range('E50').Value = '  +9.04%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '103.25'
$ws.Range('E51').Value = '  -0.97%  '
